# Tutorial 6 solution update: change date separators from "/" to "-" in
# column A (rows 3-21) and fix a few attendance tally values (D/E/G/H) for
# rows 3 and 6 to reflect the corrected attendance computation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New dash-separated dates for rows 3..21 (column A), same calendar dates
# as before, just "/" swapped for "-".
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Cells.Item($row, 1)
    # Force text storage so Excel doesn't reinterpret the dash-separated
    # string as a real date serial (the source file keeps these as plain
    # text values, matching the original slash-separated values), then
    # restore the default (unstyled) cell formatting so no stray
    # NumberFormat override is left behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
    $cell.Style = "Normal"
}

# Row 3 tally corrections: Total Attendance Count (D) 0 -> 1, Invalid (G) 0 -> 1
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 7).Value = 1

# Row 6 tally corrections: Total Attendance Count (D) 0 -> 1, Real (E) 0 -> 1,
# Absent (H) 1 -> 0
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 8).Value = 0
